$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 13 (shifts old rows 13-23 down to 14-24)
$ws.Rows(13).Insert()

# The inserted row 13 should have no A cell (clear inherited formatting)
$ws.Range("A13").Clear()

# Copy B/C column formatting (style) into the new row 13 cells
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Range("C10").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Update cell text content to the corrected / new values
$ws.Range("B10").Value = 'Desenvolver conceitos básicos da Estatística, com o apoio computacional, que permitam ao engenheiro trabalhar com o  fenômeno da aleatoriedade presente nos diversos campos de conhecimento da engenharia.'
$ws.Range("C10").Value = 'Desenvolver conceitos básicos da Estatística, com o apoio computacional, que permitam ao engenheiro trabalhar com o  fenômeno da aleatoriedade presente nos diversos campos de conhecimento da engenharia.'
$ws.Range("B13").Value = '4894221 - Mariana Pereira de Melo'
$ws.Range("C13").Value = '4894221 - Mariana Pereira de Melo'
$ws.Range("B14").Value = 'Estatística Descritiva, Modelos de Probabilidade, Teorema Central do Limite, Intervalos de Confiança, Testes de Hipóteses, ANOVA, Modelos de Regressão Linear.'
$ws.Range("C14").Value = 'Estatística Descritiva, Modelos de Probabilidade, Teorema Central do Limite, Intervalos de Confiança, Testes de Hipóteses, ANOVA, Modelos de Regressão Linear.'
$ws.Range("B16").Value = '1)Estatística Descritiva: População e amostra; apresentação gráfica dos dados; medidas de posição; medidas de dispersão.2)Amostragem: Amostragem aleatória simples com reposição; amostragem aleatória simples sem reposição.3)Conceitos de Probabilidade: Conceitos básicos de probabilidade; operações com eventos; probabilidade condicional; independência; Teorema de Bayes.4)Variáveis Aleatórias discretas: Caracterização de uma variável aleatória discreta; distribuições de probabilidade: Uniforme, Bernoulli, Binomial, Poisson, Geométrica, Binomial Negativa e Hipergeométrica.5)Variáveis Aleatórias contínuas: Caracterização de uma variável aleatória contínua; distribuições de probabilidade: Uniforme, Exponencial e Normal.6)Aproximações: Aproximação das distribuições Binomial e Poisson pela distribuição Normal.7)Teorema Central do Limite: Distribuição da média amostral; distribuição da proporção amostral; intervalos de confiança para a média amostral e para a proporção amostral; dimensionamento amostral.8)Conceitos de Testes de Hipóteses: Erro Tipo I e Erro Tipo II; p-valor; poder.9)Testes de Hipóteses para uma única amostra: Teste de hipótese para a média; teste de hipótese para a proporção e teste de hipótese para a variância.10)Testes de Hipóteses para comparação de duas amostras: Teste de hipótese para comparação de médias (amostras independentes e dependentes); teste de hipótese para comparação de duas proporções e teste de hipótese para comparação de variâncias.11) Análise de Variância: Estimação do modelo; tabela de análise de variância; intervalos de confiança para a diferença entre as médias; correção de Bonferroni; teste de homocedasticidade.12)Regressão Linear Simples e Regressão Linear Múltipla: Estimação do modelo; interpretação dos parâmetros; tabela de análise de variância; intervalos de confiança para os parâmetros; R^2; análise dos resíduos.'
$ws.Range("C16").Value = '1)Estatística Descritiva: População e amostra; apresentação gráfica dos dados; medidas de posição; medidas de dispersão.2)Amostragem: Amostragem aleatória simples com reposição; amostragem aleatória simples sem reposição.3)Conceitos de Probabilidade: Conceitos básicos de probabilidade; operações com eventos; probabilidade condicional; independência; Teorema de Bayes.4)Variáveis Aleatórias discretas: Caracterização de uma variável aleatória discreta; distribuições de probabilidade: Uniforme, Bernoulli, Binomial, Poisson, Geométrica, Binomial Negativa e Hipergeométrica.5)Variáveis Aleatórias contínuas: Caracterização de uma variável aleatória contínua; distribuições de probabilidade: Uniforme, Exponencial e Normal.6)Aproximações: Aproximação das distribuições Binomial e Poisson pela distribuição Normal.7)Teorema Central do Limite: Distribuição da média amostral; distribuição da proporção amostral; intervalos de confiança para a média amostral e para a proporção amostral; dimensionamento amostral.8)Conceitos de Testes de Hipóteses: Erro Tipo I e Erro Tipo II; p-valor; poder.9)Testes de Hipóteses para uma única amostra: Teste de hipótese para a média; teste de hipótese para a proporção e teste de hipótese para a variância.10)Testes de Hipóteses para comparação de duas amostras: Teste de hipótese para comparação de médias (amostras independentes e dependentes); teste de hipótese para comparação de duas proporções e teste de hipótese para comparação de variâncias.11) Análise de Variância: Estimação do modelo; tabela de análise de variância; intervalos de confiança para a diferença entre as médias; correção de Bonferroni; teste de homocedasticidade.12)Regressão Linear Simples e Regressão Linear Múltipla: Estimação do modelo; interpretação dos parâmetros; tabela de análise de variância; intervalos de confiança para os parâmetros; R^2; análise dos resíduos.'
$ws.Range("B19").Value = 'NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'
$ws.Range("C19").Value = 'NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'
$ws.Range("B20").Value = 'NF≥ 5,0.'
$ws.Range("C20").Value = 'NF≥ 5,0.'
$ws.Range("B21").Value = '(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.'
$ws.Range("C21").Value = '(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.'
$ws.Range("B22").Value = 'BUSSAB, Wilton O., MORETTIN, Pedro A. Estatística básica. 5. Ed. São Paulo: Saraiva, 2006.
DEVORE, Jay L Probabilidade e estatística para engenharia. São Paulo: Ed Thomson Pioneira, 2006.
JOHNSON, Richard A.; WICHERN, Dean W. Applied multivariate statistical analysis. 5. ed. Upper- Saddle River: Prentice Hall, 2002.
LARSON, Ron ; FARBER, Betsy.  Estatística aplicada. São Paulo. Ed. Prentice Hall Brasil, 2010.
HOFFMANN, R. Estatística para economistas. 4. ed. São Paulo: Pioneira, 2006.
RYAN, Thomas. Estatística moderna para engenharia. São Paulo: Ed. Campus, 2009.
RUNGER, George C.; MONTGOMERY, Douglas. Estatística aplicada e probabilidade para engenheiros. São Paulo: Ed. LTC, 2009.'
$ws.Range("C22").Value = 'BUSSAB, Wilton O., MORETTIN, Pedro A. Estatística básica. 5. Ed. São Paulo: Saraiva, 2006.
DEVORE, Jay L Probabilidade e estatística para engenharia. São Paulo: Ed Thomson Pioneira, 2006.
JOHNSON, Richard A.; WICHERN, Dean W. Applied multivariate statistical analysis. 5. ed. Upper- Saddle River: Prentice Hall, 2002.
LARSON, Ron ; FARBER, Betsy.  Estatística aplicada. São Paulo. Ed. Prentice Hall Brasil, 2010.
HOFFMANN, R. Estatística para economistas. 4. ed. São Paulo: Pioneira, 2006.
RYAN, Thomas. Estatística moderna para engenharia. São Paulo: Ed. Campus, 2009.
RUNGER, George C.; MONTGOMERY, Douglas. Estatística aplicada e probabilidade para engenheiros. São Paulo: Ed. LTC, 2009.'
